# Update the "register" sheet's selection before switching the active
# sheet, so the final active tab ends up on the newly added "search" sheet.
$wb = $excel.ActiveWorkbook

$wsRegister = $wb.Worksheets.Item("register")
$wsRegister.Range("A3:A4").Select()

# Add a new worksheet named "search" after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsSearch = $wb.Worksheets.Add($null, $lastSheet)
$wsSearch.Name = "search"

# Populate the new sheet with the search term data.
$wsSearch.Range("A1").Value = "HP"
$wsSearch.Range("A2").Value = "HP123"
$wsSearch.Range("A2").Select()
